$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E rows 2-24 contain the "purpose" value "fullRNASEQ" which should be
# corrected to "fullRNASeq" (fix formatting/casing on the fastq purpose column).
for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = "fullRNASeq"
}
